$d = $word.ActiveDocument

# --- Step 1: protect the existing run boundary right after "...many " so
# that editing the preceding run doesn't cause it to be silently re-merged
# with the following "common folk believed the " run (which carries its own
# w:rsidR and must stay a distinct run). We do this with a temporary
# "fence" bookmark, added before any other edits so the text offsets below
# are still valid against the untouched document.
$full = $d.Content.Text
$fenceMarker = "many common folk believed the"
$fencePos = $full.IndexOf($fenceMarker) + "many ".Length
$fenceRange = $d.Range($fencePos, $fencePos)
$d.Bookmarks.Add("_Fence1", $fenceRange)

# --- Step 2: insert the new "(to Eleanor Densen) " text right before
# "in what is left...".
$full = $d.Content.Text
$marker = "in what is left of the once glorious Densen family"
$splitPos = $full.IndexOf($marker)
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertBefore("(to Eleanor Densen) ")

# --- Step 3: split "(to Eleanor Densen) " into its own run (distinct from
# both " You are a cousin " and "in what is left...") and re-home the
# "_GoBack" bookmark directly after it. Adding a bookmark forces a run
# boundary at that character position; adding one named "_GoBack" also
# automatically removes/replaces any pre-existing "_GoBack" bookmark
# elsewhere in the document (Word keeps only one).
$full = $d.Content.Text
$phrase = "(to Eleanor Densen) "
$idx1 = $full.IndexOf($phrase)
$idx2 = $idx1 + $phrase.Length

$beforeRange = $d.Range($idx1, $idx1)
$d.Bookmarks.Add("_TempSplit", $beforeRange)

$afterRange = $d.Range($idx2, $idx2)
$d.Bookmarks.Add("_GoBack", $afterRange)

# --- Step 4: clean up the temporary bookmarks. Their run-splitting effect
# persists even after the bookmarks themselves are removed.
$d.Bookmarks("_TempSplit").Delete()
$d.Bookmarks("_Fence1").Delete()

Write-Output "done"
